$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '42.095.58'
$r.ClearFormats()
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '2.261.58'
$r.ClearFormats()
$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  +2.34%  '
$r.ClearFormats()
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  -0.04%  '
$r.ClearFormats()
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '302.52'
$r.ClearFormats()
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +3.80%  '
$r.ClearFormats()
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '92.92'
$r.ClearFormats()
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  +6.97%  '
$r.ClearFormats()
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  +3.73%  '
$r.ClearFormats()
$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  -0.04%  '
$r.ClearFormats()
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  +3.74%  '
$r.ClearFormats()
$r = $ws.Range('B10')
$r.NumberFormat = '@'
$r.Value = 'OKB'
$r.ClearFormats()
$r = $ws.Range('C10')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$r.ClearFormats()
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '54.77'
$r.ClearFormats()
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  +9.53%  '
$r.ClearFormats()
$r = $ws.Range('B11')
$r.NumberFormat = '@'
$r.Value = 'Avalanche'
$r.ClearFormats()
$r = $ws.Range('C11')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$r.ClearFormats()
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '32.69'
$r.ClearFormats()
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  +8.14%  '
$r.ClearFormats()
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  +2.64%  '
$r.ClearFormats()
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  +4.09%  '
$r.ClearFormats()
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '2.613.32'
$r.ClearFormats()
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  +2.51%  '
$r.ClearFormats()
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '14.15'
$r.ClearFormats()
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +3.13%  '
$r.ClearFormats()
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '2.267.97'
$r.ClearFormats()
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  +3.76%  '
$r.ClearFormats()
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  +3.85%  '
$r.ClearFormats()
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '41.954.03'
$r.ClearFormats()
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  +5.49%  '
$r.ClearFormats()
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '12.22'
$r.ClearFormats()
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +9.49%  '
$r.ClearFormats()
$r = $ws.Range('D21')
$r.Value = '0.0Q0909'
$r.Replace('Q', [char]0x2083) | Out-Null
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  +3.72%  '
$r.ClearFormats()
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '67.33'
$r.ClearFormats()
$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  +2.81%  '
$r.ClearFormats()
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '241.87'
$r.ClearFormats()
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  +2.04%  '
$r.ClearFormats()
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  +5.81%  '
$r.ClearFormats()
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  -0.12%  '
$r.ClearFormats()
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +4.60%  '
$r.ClearFormats()
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '23.99'
$r.ClearFormats()
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  +3.78%  '
$r.ClearFormats()
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  +7.13%  '
$r.ClearFormats()
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '9.70'
$r.ClearFormats()
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  +5.35%  '
$r.ClearFormats()
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '34.25'
$r.ClearFormats()
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +7.54%  '
$r.ClearFormats()
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '159.30'
$r.ClearFormats()
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +1.31%  '
$r.ClearFormats()
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  +0.00%  '
$r.ClearFormats()
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '5.17'
$r.ClearFormats()
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  +4.56%  '
$r.ClearFormats()
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '0.0744'
$r.ClearFormats()
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  +5.10%  '
$r.ClearFormats()
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '3.06'
$r.ClearFormats()
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  +4.42%  '
$r.ClearFormats()
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  +3.06%  '
$r.ClearFormats()
$r = $ws.Range('B38')
$r.NumberFormat = '@'
$r.Value = 'Celestia'
$r.ClearFormats()
$r = $ws.Range('C38')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$r.ClearFormats()
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '16.66'
$r.ClearFormats()
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  +9.52%  '
$r.ClearFormats()
$r = $ws.Range('B39')
$r.NumberFormat = '@'
$r.Value = 'Kaspa'
$r.ClearFormats()
$r = $ws.Range('C39')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$r.ClearFormats()
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.104'
$r.ClearFormats()
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  +6.57%  '
$r.ClearFormats()
$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  +4.65%  '
$r.ClearFormats()
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  +5.10%  '
$r.ClearFormats()
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  +6.46%  '
$r.ClearFormats()
$r = $ws.Range('B43')
$r.NumberFormat = '@'
$r.Value = 'EnergySwap'
$r.ClearFormats()
$r = $ws.Range('C43')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r.ClearFormats()
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '20.18'
$r.ClearFormats()
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +13.39%  '
$r.ClearFormats()
$r = $ws.Range('B44')
$r.NumberFormat = '@'
$r.Value = 'Maker'
$r.ClearFormats()
$r = $ws.Range('C44')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$r.ClearFormats()
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '2.055.48'
$r.ClearFormats()
$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  -2.64%  '
$r.ClearFormats()
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  +4.00%  '
$r.ClearFormats()
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '10.13'
$r.ClearFormats()
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  +1.17%  '
$r.ClearFormats()
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '2.91'
$r.ClearFormats()
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  +8.39%  '
$r.ClearFormats()
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  -4.28%  '
$r.ClearFormats()
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '2.484.08'
$r.ClearFormats()
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  +2.56%  '
$r.ClearFormats()
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '1.53'
$r.ClearFormats()
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  +4.32%  '
$r.ClearFormats()
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  +4.87%  '
$r.ClearFormats()
